# Commit all constraints feasible model
# Update the binary decision-variable matrices on sheets y1, y2, y3.

$wb = $excel.ActiveWorkbook

# --- Sheet y1 ---
$ws1 = $wb.Worksheets.Item("y1")
$ws1.Range("A3").Value = 0
$ws1.Range("A4").Value = 1
$ws1.Range("A5").Value = 1
$ws1.Range("A7").Value = 1
$ws1.Range("A8").Value = 0

# --- Sheet y2 ---
$ws2 = $wb.Worksheets.Item("y2")
$ws2.Range("A6").Value = 0

# --- Sheet y3 ---
$ws3 = $wb.Worksheets.Item("y3")
$ws3.Range("A3").Value = 1
$ws3.Range("A4").Value = 0
$ws3.Range("A5").Value = 0
$ws3.Range("B6").Value = 1
$ws3.Range("A7").Value = 0
$ws3.Range("A8").Value = 1
